$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create new last row (194) as a copy of row 193 (which will become the
# "new" last data point once everything below row 22 shifts down by one).
$lastDataCol = 18
for ($c = 1; $c -le $lastDataCol; $c++) {
    $ws.Cells.Item(194, $c).Value2 = $ws.Cells.Item(193, $c).Value2
}
$ws.Cells.Item(194, 4).NumberFormat = $ws.Cells.Item(193, 4).NumberFormat

# Shift the Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) columns down by one
# row, for every row from 193 down to 23 (each row takes on the values
# that used to belong to the row directly above it). Walking top-to-bottom
# in descending row order means every source row is read before it gets
# overwritten.
for ($r = 193; $r -ge 23; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value2 = $ws.Cells.Item($src, 4).Value2
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($src, 10).Value2
    $ws.Cells.Item($r, 11).Value2 = $ws.Cells.Item($src, 11).Value2
    $ws.Cells.Item($r, 12).Value2 = $ws.Cells.Item($src, 12).Value2
    $ws.Cells.Item($r, 13).Value2 = $ws.Cells.Item($src, 13).Value2
    $ws.Cells.Item($r, 16).Value2 = $ws.Cells.Item($src, 16).Value2
}

# Row 22 is the new top of the shifted block; it gets a brand new data
# point that did not exist before (Fecha 2021-10-21 / Volumen 160).
$ws.Cells.Item(22, 4).Value2 = 44490
$ws.Cells.Item(22, 10).Value2 = 160
